# "Advanced Excel - For Analysis": add a second worksheet "решение" that
# holds the Solver-produced solution for the budget-allocation model that
# already lives on "Лист1", wire up the Excel Solver parameters (stored as
# hidden, sheet-scoped defined names) for that new sheet, and make the new
# sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Duplicate "Лист1" right after itself -> becomes the new sheet 2.
# ------------------------------------------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "решение"

# ------------------------------------------------------------------
# 2. The copy/rename round-trip can mangle the handful of pre-existing
#    "Лист1!#REF!" hidden defined names (they lose their sheet prefix).
#    Put them back exactly as they were before we touched anything else.
# ------------------------------------------------------------------
$nameCount = $ws1.Names.Count
for ($i = 1; $i -le $nameCount; $i++) {
    $nm = $ws1.Names.Item($i)
    if ($nm.RefersTo -eq "=#REF!") {
        $nm.RefersTo = "=Лист1!#REF!"
    }
}

# ------------------------------------------------------------------
# 3. Fill in the Solver's adjustable cells on the new sheet with the
#    values Solver converged on. The dependent formulas already copied
#    over from "Лист1" (row 7 totals, row 8 ratios, row 9 spend, etc.)
#    recompute automatically from these inputs.
# ------------------------------------------------------------------
$ws2.Range("B7").Value = 2750
$ws2.Range("C7").Value = 5000
$ws2.Range("D7").Value = 0
$ws2.Range("E7").Value = 3600

# ------------------------------------------------------------------
# 4. Match the selection/active-tab state: both sheets end up with B11
#    selected, and "решение" becomes the active (visible) tab.
# ------------------------------------------------------------------
$ws1.Range("B11").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B11").Select() | Out-Null

# ------------------------------------------------------------------
# 5. Recreate the Excel Solver bookkeeping (hidden, sheet-scoped
#    defined names) for the new "решение" sheet, mirroring the set
#    already present for "Лист1".
# ------------------------------------------------------------------
$solverNames = @(
    @{Name="solver_adj";  Value="решение!`$B`$7:`$E`$7"},
    @{Name="solver_cvg";  Value="0.0001"},
    @{Name="solver_drv";  Value="1"},
    @{Name="solver_eng";  Value="2"},
    @{Name="solver_est";  Value="1"},
    @{Name="solver_itr";  Value="2147483647"},
    @{Name="solver_lhs1"; Value="решение!`$B`$7:`$E`$7"},
    @{Name="solver_lhs2"; Value="решение!`$B`$8:`$E`$8"},
    @{Name="solver_lhs3"; Value="решение!`$F`$9"},
    @{Name="solver_lhs4"; Value="решение!#REF!"},
    @{Name="solver_mip";  Value="2147483647"},
    @{Name="solver_mni";  Value="30"},
    @{Name="solver_mrt";  Value="0.075"},
    @{Name="solver_msl";  Value="2"},
    @{Name="solver_neg";  Value="1"},
    @{Name="solver_nod";  Value="2147483647"},
    @{Name="solver_num";  Value="3"},
    @{Name="solver_nwt";  Value="1"},
    @{Name="solver_opt";  Value="решение!`$F`$7"},
    @{Name="solver_pre";  Value="0.000001"},
    @{Name="solver_rbv";  Value="2"},
    @{Name="solver_rel1"; Value="3"},
    @{Name="solver_rel2"; Value="1"},
    @{Name="solver_rel3"; Value="2"},
    @{Name="solver_rel4"; Value="2"},
    @{Name="solver_rhs1"; Value="0"},
    @{Name="solver_rhs2"; Value="решение!`$B`$4:`$E`$4"},
    @{Name="solver_rhs3"; Value="решение!`$B`$11"},
    @{Name="solver_rhs4"; Value="решение!#REF!"},
    @{Name="solver_rlx";  Value="2"},
    @{Name="solver_rsd";  Value="0"},
    @{Name="solver_scl";  Value="2"},
    @{Name="solver_sho";  Value="2"},
    @{Name="solver_ssz";  Value="0"},
    @{Name="solver_tim";  Value="2147483647"},
    @{Name="solver_tol";  Value="0.01"},
    @{Name="solver_typ";  Value="2"},
    @{Name="solver_val";  Value="0"},
    @{Name="solver_ver";  Value="3"}
)

foreach ($item in $solverNames) {
    $ws2.Names.Add($item.Name, "=" + $item.Value) | Out-Null
}

# Mark every one of those as hidden (matches Solver's own bookkeeping).
# Must be done by fully-qualified name lookup - indexed/ambiguous access
# to Names can resolve to the same-named entry on "Лист1" instead.
foreach ($item in $solverNames) {
    $qualified = "решение!" + $item.Name
    $nm = $ws2.Names.Item($qualified)
    $nm.Visible = $false
}

$wb.Save()
